$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 547; this shifts the existing rows 547-600 down
# to 548-601 (row 600's old data lands on 601, matching the new tail row).
$ws.Rows.Item(547).Insert()

# Populate the newly inserted row 547 with the new price record.
$ws.Cells.Item(547, 1).Value = 9
$ws.Cells.Item(547, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(547, 3).Value = "Metropolitana"
$ws.Cells.Item(547, 4).Value = 45106
$ws.Cells.Item(547, 5).Value = 13
$ws.Cells.Item(547, 6).Value = 100112052
$ws.Cells.Item(547, 7).Value = "Albahaca"
$ws.Cells.Item(547, 8).Value = "Sin especificar"
$ws.Cells.Item(547, 9).Value = "Primera"
$ws.Cells.Item(547, 10).Value = 250
$ws.Cells.Item(547, 11).Value = 4200
$ws.Cells.Item(547, 12).Value = 4500
$ws.Cells.Item(547, 13).Value = 4350
$ws.Cells.Item(547, 14).Value = "$/paquete"
$ws.Cells.Item(547, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(547, 16).Value = 4350
$ws.Cells.Item(547, 17).Value = 1
$ws.Cells.Item(547, 18).Value = "Hortaliza"
